$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells retain their original text formatting (avoid Excel
# auto-converting numeric-looking strings like "63.548.17" or
# "0.0000237" into actual numbers / scientific notation).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.548.17'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.065.56'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.38'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.08'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.064.86'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.156'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.85'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.85'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000237'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.571.71'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.20'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.501.30'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.070.04'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '490.95'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.44'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.707'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.46'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.88'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.87'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.70'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +10.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.42'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.69'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.21'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.30'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0820'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.31'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.99'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.27'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.64'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '437.84'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.56%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0363'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.841.50'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.24'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.27'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.34'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.56%  '
